# repull data, push all data, mean calculation
# Update column F (dSF) values that changed after re-pulling the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -11
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = -3
$ws.Range("F22").Value = 0
